$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 996, shifting existing rows 996:1072 down to 997:1073
$ws.Rows(996).Insert()

# Populate the newly inserted row 996 with the new data record
$ws.Range("A996").Value = 5
$ws.Range("B996").Value = "Macroferia Regional de Talca"
$ws.Range("C996").Value = "Maule"
$ws.Range("D996").Value = 45265
$ws.Range("E996").Value = 7
$ws.Range("F996").Value = "Fruta"
$ws.Range("G996").Value = 100102
$ws.Range("H996").Value = "Cítricos"
$ws.Range("I996").Value = 100102005
$ws.Range("J996").Value = "Naranja"
$ws.Range("K996").Value = "Valencia"
$ws.Range("L996").Value = "Primera"
$ws.Range("M996").Value = 350
$ws.Range("N996").Value = 11000
$ws.Range("O996").Value = 11000
$ws.Range("P996").Value = 11000
$ws.Range("Q996").Value = "$/bandeja 15 kilos granel"
$ws.Range("R996").Value = "Provincia de Curicó"
$ws.Range("S996").Value = 733
$ws.Range("T996").Value = 15
